# Fixed empty input start week case
# This chronogram/Gantt template generator previously produced actual
# month/date-range labels (e.g. "December" / "25/Dec - 31/Dec") computed
# from a start-date input. When that input is empty, the generator now
# falls back to generic placeholders ("Month 1" / "Week 1", ...), and the
# trailing, not-fully-populated month column is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Break apart the merges that are going to be restructured. Leave the
#    row-spanning merges for the B/C/D/E header columns untouched.
# ---------------------------------------------------------------------
$ws.Range("G1").UnMerge()
$ws.Range("F1").UnMerge()
$ws.Range("F2").UnMerge()
$ws.Range("G2").UnMerge()
$ws.Range("L2").UnMerge()
$ws.Range("P2").UnMerge()

# ---------------------------------------------------------------------
# 2. Row 1: drop the second year header (G1, "2025"). F1 keeps "2024".
# ---------------------------------------------------------------------
$ws.Range("G1").ClearContents()

# ---------------------------------------------------------------------
# 3. Row 2: generic month placeholders instead of real month names.
#    J2/N2 previously held no data at all, so pick up the shared
#    "month band" look (the same style used by F2/G2/L2/P2) first.
# ---------------------------------------------------------------------
$ws.Range("F2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("N2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("F2").Value = "Month 1"
$ws.Range("G2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("J2").Value = "Month 2"
$ws.Range("N2").Value = "Month 3"

# ---------------------------------------------------------------------
# 4. Row 3: generic week placeholders instead of real date ranges.
# ---------------------------------------------------------------------
$ws.Range("F3").Value = "Week 1"
$ws.Range("G3").Value = "Week 2"
$ws.Range("H3").Value = "Week 3"
$ws.Range("I3").Value = "Week 4"
$ws.Range("J3").Value = "Week 5"
$ws.Range("K3").Value = "Week 6"
$ws.Range("L3").Value = "Week 7"
$ws.Range("M3").Value = "Week 8"
$ws.Range("N3").Value = "Week 9"
$ws.Range("O3").Value = "Week 10"
$ws.Range("P3").Value = "Week 11"
$ws.Range("Q3").ClearContents()

# ---------------------------------------------------------------------
# 5. Rows 4-13: the explicit Start/End date helper columns (D, E) are no
#    longer populated when the start date input is empty.
# ---------------------------------------------------------------------
$ws.Range("D4:E13").ClearContents()

# ---------------------------------------------------------------------
# 6. Give the trailing, incomplete month cell (P2) the "header" look
#    (bold white text on the blue fill) without the centring alignment
#    used elsewhere, then drop column Q altogether (it only ever held
#    the now-removed 12th week / 2025 header).
# ---------------------------------------------------------------------
$ws.Range("P2").ClearContents()
$ws.Range("F1").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").HorizontalAlignment = 1
$ws.Range("P2").VerticalAlignment = -4107
$ws.Application.CutCopyMode = $false

$ws.Columns("Q").Delete()

# ---------------------------------------------------------------------
# 7. Re-create the merges with the new month/week groupings.
# ---------------------------------------------------------------------
$ws.Range("F1:P1").Merge()
$ws.Range("F2:I2").Merge()
$ws.Range("J2:M2").Merge()
$ws.Range("N2:O2").Merge()

Write-Output "done"
